$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C3").Value = 9963
$ws.Range("C4:C5").Value = 9938
$ws.Range("C6").Value = 9071
$ws.Range("C7:C9").Value = 9019
$ws.Range("C10:C11").Value = 8976
$ws.Range("C12:C13").Value = 8793
$ws.Range("C14:C30").Value = 8733
$ws.Range("C31:C33").Value = 8716
$ws.Range("C34").Value = 8490
$ws.Range("C35:C46").Value = 8085
$ws.Range("C56").Value = 8013
$ws.Range("C57:C77").Value = 7590
$ws.Range("C102:C153").Value = 7573
